$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'63.180.28"
$ws.Range('E2').Value = '  -3.16%  '
$ws.Range('D3').Value = "'3.080.62"
$ws.Range('E3').Value = '  -2.10%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = "'547.46"
$ws.Range('E5').Value = '  -3.31%  '
$ws.Range('D6').Value = "'136.31"
$ws.Range('E6').Value = '  -8.22%  '
$ws.Range('E7').Value = '  +0.17%  '
$ws.Range('D8').Value = "'3.070.00"
$ws.Range('E8').Value = '  -2.29%  '
$ws.Range('D9').Value = "'0.492"
$ws.Range('E9').Value = '  -1.18%  '
$ws.Range('D10').Value = "'6.57"
$ws.Range('E10').Value = '  -5.21%  '
$ws.Range('D11').Value = "'0.157"
$ws.Range('E11').Value = '  -0.55%  '
$ws.Range('D12').Value = "'0.460"
$ws.Range('E12').Value = '  -0.53%  '
$ws.Range('D13').Value = "'34.93"
$ws.Range('E13').Value = '  -3.13%  '
$ws.Range('D14').Value = "'0.0000216"
$ws.Range('E14').Value = '  -2.44%  '
$ws.Range('D15').Value = "'3.579.85"
$ws.Range('E15').Value = '  -1.96%  '
$ws.Range('D16').Value = "'63.329.29"
$ws.Range('E16').Value = '  -3.02%  '
$ws.Range('E17').Value = '  -1.21%  '
$ws.Range('D18').Value = "'3.086.63"
$ws.Range('E18').Value = '  -1.76%  '
$ws.Range('D19').Value = "'6.64"
$ws.Range('E19').Value = '  -1.49%  '
$ws.Range('D20').Value = "'482.17"
$ws.Range('E20').Value = '  -8.23%  '
$ws.Range('D21').Value = "'13.39"
$ws.Range('E21').Value = '  -3.28%  '
$ws.Range('D22').Value = "'0.700"
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('D23').Value = "'7.15"
$ws.Range('E23').Value = '  -3.58%  '
$ws.Range('D24').Value = "'77.63"
$ws.Range('E24').Value = '  -1.35%  '
$ws.Range('D25').Value = "'12.18"
$ws.Range('E25').Value = '  -4.31%  '
$ws.Range('E26').Value = '  +0.08%  '
$ws.Range('D27').Value = "'2.73"
$ws.Range('E27').Value = '  -2.66%  '
$ws.Range('D28').Value = "'8.24"
$ws.Range('E28').Value = '  -4.68%  '
$ws.Range('E29').Value = '  +0.07%  '
$ws.Range('D30').Value = "'1.93"
$ws.Range('E30').Value = '  -9.25%  '
$ws.Range('D31').Value = "'26.32"
$ws.Range('E31').Value = '  +0.44%  '
$ws.Range('D32').Value = "'1.13"
$ws.Range('E32').Value = '  -0.55%  '
$ws.Range('B33').Value = 'OKB'
$ws.Range('C33').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D33').Value = "'61.42"
$ws.Range('E33').Value = '  +15.65%  '
$ws.Range('B34').Value = 'Stacks'
$ws.Range('C34').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D34').Value = "'2.49"
$ws.Range('E34').Value = '  -6.94%  '
$ws.Range('D35').Value = "'530.57"
$ws.Range('E35').Value = '  -4.66%  '
$ws.Range('D36').Value = "'5.93"
$ws.Range('E36').Value = '  -2.21%  '
$ws.Range('D37').Value = "'5.15"
$ws.Range('E37').Value = '  -5.37%  '
$ws.Range('D38').Value = "'0.0400"
$ws.Range('E38').Value = '  -9.67%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = "'3.084.72"
$ws.Range('E39').Value = '  +0.19%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').Value = "'0.0789"
$ws.Range('E40').Value = '  -3.99%  '
$ws.Range('D41').Value = "'0.118"
$ws.Range('E41').Value = '  -2.62%  '
$ws.Range('B42').Value = 'Cosmos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D42').Value = "'8.08"
$ws.Range('E42').Value = '  -2.16%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').Value = "'2.66"
$ws.Range('E43').Value = '  -7.65%  '
$ws.Range('D44').Value = "'0.253"
$ws.Range('E44').Value = '  -1.32%  '
$ws.Range('E46').Value = '  -6.87%  '
$ws.Range('D47').Value = "'121.63"
$ws.Range('E47').Value = '  +2.82%  '
$ws.Range('D48').Value = "'24.24"
$ws.Range('E48').Value = '  -2.94%  '
$ws.Range('E49').Value = '  -1.82%  '
$ws.Range('D50').Value = "'0.0₃0503"
$ws.Range('E50').Value = '  -4.28%  '
$ws.Range('D51').Value = "'2.34"
$ws.Range('E51').Value = '  +60.18%  '
